# Agrega soporte para TABLE
# - Quita la columna de demo "Merge" (E2:F5, con sus celdas combinadas)
# - Completa el encabezado "Dato" en A15 (tabla de maquinas)
# - Agrega la tabla XOR en A20:C22
# - Registra los nombres PUE.TABLE.tabla_maquinas y PUE.TABLE.tabla_xor

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Eliminar la columna "Merge" de ejemplo (E2:F5), incluidas las celdas combinadas
$ws.Range("E2:F5").UnMerge()
$ws.Range("E2:F5").Clear()

# 2) La tabla "tabla_maquinas" (A15:D18) ya existe; solo falta el rotulo de la
#    primera columna del encabezado
$ws.Range("A15").Value = "Dato"

# 3) Nueva tabla "tabla_xor" (A20:C22), con el mismo estilo que la tabla de
#    maquinas: fila de encabezado como A15, filas de datos como las de A16:D18
$ws.Range("A15").Copy($ws.Range("A20"))
$ws.Range("B16:C16").Copy($ws.Range("B20:C20"))
$ws.Range("A16:C16").Copy($ws.Range("A21:C21"))
$ws.Range("A17:C17").Copy($ws.Range("A22:C22"))

$ws.Range("A20").Value = "XOR"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 1
$ws.Range("A21").Value = 0
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 1
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 0

# 4) Nombres definidos para las tablas
$wb.Names.Add("PUE.TABLE.tabla_maquinas", "=Sheet1!`$A`$15:`$D`$18")
$wb.Names.Add("PUE.TABLE.tabla_xor", "=Sheet1!`$A`$20:`$C`$22")

# 5) Selección final tal como quedó en el archivo editado
$ws.Range("C18").Select()
